$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared text value used in column A (A2:A5 all share this string)
# so that "pizza+delivery" becomes "delivery+person"
$ws.Range("A2:A5").Value = "delivery+person"

# Widen column A (20.2 "characters" serializes to a stored width of 21)
$ws.Range("A1").EntireColumn.ColumnWidth = 20.2

# Move the active cell selection to A6
$ws.Range("A6").Select()

$wb.Save()
